$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# Consolidate translation columns C/D into column B:
# prefer D's value, fall back to C's value, fall back to A's value
# (mirrors: D non-empty -> use D; elseif C non-empty -> use C; else -> use A)
for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)

    $cVal = $cCell.Value2
    $dVal = $dCell.Value2

    if ($dVal -ne $null -and $dVal -ne "") {
        $dCell.Copy($ws.Cells.Item($r, 2))
    } elseif ($cVal -ne $null -and $cVal -ne "") {
        $cCell.Copy($ws.Cells.Item($r, 2))
    } else {
        $ws.Cells.Item($r, 1).Copy($ws.Cells.Item($r, 2))
    }
}

for ($r = 1; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).ClearContents()
    $ws.Cells.Item($r, 4).ClearContents()
}
